$d = $word.ActiveDocument

# 1. Merge "Illustrat" + bookmark + "ion/Paint" into a single run "Illustration/Paint".
#    Word's Find/Replace treats the bookmark as a zero-width boundary in the text
#    stream, so searching for the already-combined text and replacing it in place
#    collapses the two runs (and drops the now-redundant bookmark) into one run.
$null = $d.Content.Find.Execute("Illustration/Paint", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "Illustration/Paint", 2)

# 2. Insert a new paragraph after "Logos" for the bold "Color Pallette:" line,
#    including the spell-check markers Word leaves around "Pallette".
$logosIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq "Logos") {
        $logosIndex = $i
    }
}
$logosPara = $d.Paragraphs.Item($logosIndex)
$logosPara.Range.InsertParagraphAfter()

$colorPara = $d.Paragraphs.Item($logosIndex + 1)
$colorXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
  '<w:p>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Color </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Pallette</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r>' + `
  '</w:p>' + `
  '</w:body></w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'
$null = $colorPara.Range.InsertXML($colorXml)

# 3. Insert a new paragraph after that for the bold "Paragraph about yourself:" line.
$colorPara = $d.Paragraphs.Item($logosIndex + 1)
$colorPara.Range.InsertParagraphAfter()

$aboutPara = $d.Paragraphs.Item($logosIndex + 2)
$aboutXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
  '<w:p>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Paragraph about yourself:</w:t></w:r>' + `
  '</w:p>' + `
  '</w:body></w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'
$null = $aboutPara.Range.InsertXML($aboutXml)

# 4. Re-home the "_GoBack" bookmark onto the trailing empty paragraph.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$null = $d.Bookmarks.Add("_GoBack", $lastPara.Range)
